$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing GNN-MT / PN / RF rows shift down one slot to make room for the
# newly-reported "GNN-MT-O" series, and every mean-rank value is refreshed
# now that the ranking includes more models.
$ws.Range("A2").Value = "16_train (GNN-MT-O) val delta-auprc"
$ws.Range("B2").Value = 3.738095238095238

$ws.Range("A3").Value = "16_train (GNN-MT) val delta-auprc"
$ws.Range("B3").Value = 3.095238095238095

$ws.Range("A4").Value = "16_train (PN) val delta-auprc"
$ws.Range("B4").Value = 2.952380952380953

$ws.Range("A5").Value = "16_train (RF) val delta-auprc"
$ws.Range("B5").Value = 2.642857142857143

$ws.Range("A6").Value = "16_train (PN-O) val delta-auprc"
$ws.Range("B6").Value = 2.571428571428572

# Copy the label formatting (bold, bordered, centered/top-aligned) down
# onto the two newly populated label cells.
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A6").PasteSpecial(-4122)   # xlPasteFormats
